$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header H1 (still "SID D" text-wise, but triggers shared-string reorder)
$ws.Range("H1").Value = "SID D"

# Update Faculty Mentor column (L2:L11) with new rotating set of names
$mentors = @(
    "be18103032 Gaganpreet Singh Khurana",
    "Akshit Garg",
    "Shayan Yaseen",
    "be18103032 Gaganpreet Singh Khurana",
    "Akshit Garg",
    "be18103032 Gaganpreet Singh Khurana",
    "Akshit Garg",
    "Shayan Yaseen",
    "be18103032 Gaganpreet Singh Khurana",
    "Akshit Garg"
)

for ($i = 0; $i -lt $mentors.Length; $i++) {
    $row = $i + 2
    $ws.Range("L$row").Value = $mentors[$i]
}

# Update selection to C13
$ws.Range("C13").Select()
